$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.685.96'
$ws.Range('D3').Value = '3.554.79'
$ws.Range('E3').Value = '  +1.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.29'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '187.04'
$ws.Range('E6').Value = '  +1.51%  '
$ws.Range('E7').Value = '  +2.30%  '
$ws.Range('D8').Value = '3.544.83'
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.221'
$ws.Range('E10').Value = '  +18.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.648'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.74'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('E13').Value = '  +6.17%  '
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '4.122.13'
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Value = '70.721.79'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.592.68'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.13'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.76'
$ws.Range('E19').Value = '  +4.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '572.36'
$ws.Range('E20').Value = '  +6.58%  '
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.68'
$ws.Range('E23').Value = '  -5.33%  '
$ws.Range('E24').Value = '  +3.17%  '
$ws.Range('E25').Value = '  -1.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '93.84'
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.24'
$ws.Range('E27').Value = '  +4.26%  '
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.27'
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.48'
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.30'
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('E33').Value = '  +1.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.04'
$ws.Range('E34').Value = '  -2.82%  '
$ws.Range('B35').Value = 'dogwifhat'
$ws.Range('C35').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.74'
$ws.Range('E35').Value = '  +19.08%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.36'
$ws.Range('E36').Value = '  +11.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '543.00'
$ws.Range('E37').Value = '  -3.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.412'
$ws.Range('E38').Value = '  +3.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.11'
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('D40').Value = '0.0₃0803'
$ws.Range('E40').Value = '  +4.87%  '
$ws.Range('D42').Value = '3.586.25'
$ws.Range('E42').Value = '  +11.19%  '
$ws.Range('E43').Value = '  +3.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.45'
$ws.Range('E44').Value = '  +3.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0471'
$ws.Range('E45').Value = '  +6.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.50'
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.34'
$ws.Range('E48').Value = '  +4.24%  '
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.47'
$ws.Range('E50').Value = '  +8.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  +0.14%  '
